# Append a new log row (row 19) to the Nalco run log, mirroring the
# style/format of the previous row (row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing row down to the new row
# so the new row's cells keep the same style index (centered, etc.).
$ws.Range("A18:H18").Copy()
$ws.Range("A19:H19").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new run-log entry.
$ws.Cells.Item(19, 1).Value = "2025-08-16 06:47:43 UTC"
$ws.Cells.Item(19, 2).Value = "2025-08-16 12:17:43 IST"
$ws.Cells.Item(19, 3).Value = "SKIPPED"
$ws.Cells.Item(19, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(19, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item(19, 7).Value = 0

Write-Host "Appended row 19 to run log."
